$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 634.9259
$ws.Range("I33").Value = 562.9091
$ws.Range("K33").Value = 562.9091
$ws.Range("M33").Value = -333.9091
$ws.Range("H132").Value = 2613.5305
$ws.Range("I132").Value = 756.02856
$ws.Range("K132").Value = 2268.08568
$ws.Range("M132").Value = 261.9143199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 28221
$ws.Range("J44").Value = 28221
$ws.Range("L44").Value = 28221
$ws.Range("N44").Value = -29197
$ws.Range("H55").Value = 16072
$ws.Range("J55").Value = 16072
$ws.Range("L55").Value = 16072
$ws.Range("N55").Value = -16702
$ws.Range("H61").Value = 3744
$ws.Range("I61").Value = 2583.3635
$ws.Range("J61").Value = 7999.6665
$ws.Range("K61").Value = 2583.3635
$ws.Range("L61").Value = 7999.6665
$ws.Range("M61").Value = -2371.3635
$ws.Range("N61").Value = -8423.666499999999
$ws.Range("H80").Value = 30855
$ws.Range("J80").Value = 30855
$ws.Range("L80").Value = 30855
$ws.Range("N80").Value = -32851
$ws.Range("H83").Value = 30855
$ws.Range("J83").Value = 30855
$ws.Range("L83").Value = 92565
$ws.Range("N83").Value = -102549
$ws.Range("H136").Value = 3744
$ws.Range("I136").Value = 2583.3635
$ws.Range("J136").Value = 7999.6665
$ws.Range("K136").Value = 7750.0905
$ws.Range("L136").Value = 23998.9995
$ws.Range("M136").Value = -5200.0905
$ws.Range("N136").Value = -29098.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 24897.555
$ws.Range("J35").Value = 24897.555
$ws.Range("L35").Value = 24897.555
$ws.Range("N35").Value = -25517.555
$ws.Range("H82").Value = 15350.4
$ws.Range("J82").Value = 32792.668
$ws.Range("L82").Value = 32792.668
$ws.Range("N82").Value = -33558.668
$ws.Range("H85").Value = 15350.4
$ws.Range("J85").Value = 32792.668
$ws.Range("L85").Value = 32792.668
$ws.Range("N85").Value = -35444.668
$ws.Range("H94").Value = 1092.9412
$ws.Range("I94").Value = 1093.8182
$ws.Range("K94").Value = 1093.8182
$ws.Range("M94").Value = -642.8181999999999
$ws.Range("H122").Value = 52271.43
$ws.Range("J122").Value = 52271.43
$ws.Range("L122").Value = 52271.43
$ws.Range("N122").Value = -62071.43
$ws.Range("H125").Value = 53904.445
$ws.Range("J125").Value = 53904.445
$ws.Range("L125").Value = 53904.445
$ws.Range("N125").Value = -63744.445

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2755.577
$ws.Range("I31").Value = 1800.6316
$ws.Range("J31").Value = 3305.394
$ws.Range("K31").Value = 1800.6316
$ws.Range("L31").Value = 3305.394
$ws.Range("M31").Value = -1505.6316
$ws.Range("N31").Value = -3895.394
$ws.Range("H34").Value = 2755.577
$ws.Range("I34").Value = 1800.6316
$ws.Range("J34").Value = 3305.394
$ws.Range("K34").Value = 1800.6316
$ws.Range("L34").Value = 3305.394
$ws.Range("M34").Value = -1598.6316
$ws.Range("N34").Value = -3709.394
$ws.Range("H41").Value = 13397.714
$ws.Range("J41").Value = 14637.5
$ws.Range("L41").Value = 14637.5
$ws.Range("N41").Value = -15493.5
$ws.Range("H50").Value = 10875.8
$ws.Range("J50").Value = 10875.8
$ws.Range("L50").Value = 10875.8
$ws.Range("N50").Value = -12125.8
$ws.Range("H51").Value = 18029.545
$ws.Range("J51").Value = 18029.545
$ws.Range("L51").Value = 18029.545
$ws.Range("N51").Value = -19501.545
$ws.Range("H60").Value = 15800.75
$ws.Range("J60").Value = 19067.666
$ws.Range("L60").Value = 19067.666
$ws.Range("N60").Value = -20089.666
$ws.Range("H61").Value = 18029.545
$ws.Range("J61").Value = 18029.545
$ws.Range("L61").Value = 18029.545
$ws.Range("N61").Value = -18725.545
$ws.Range("H68").Value = 25268.715
$ws.Range("J68").Value = 25268.715
$ws.Range("L68").Value = 25268.715
$ws.Range("N68").Value = -26766.715
$ws.Range("H71").Value = 25268.715
$ws.Range("J71").Value = 25268.715
$ws.Range("L71").Value = 75806.145
$ws.Range("N71").Value = -83294.145
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 82500
$ws.Range("J87").Value = 76666.664
$ws.Range("L87").Value = 229999.992
$ws.Range("N87").Value = -232495.992
$ws.Range("H90").Value = 82500
$ws.Range("J90").Value = 76666.664
$ws.Range("L90").Value = 689999.976
$ws.Range("N90").Value = -702479.976
$ws.Range("H107").Value = 936.7917
$ws.Range("I107").Value = 584.7646999999999
$ws.Range("J107").Value = 1791.7142
$ws.Range("K107").Value = 1754.2941
$ws.Range("L107").Value = 5375.142599999999
$ws.Range("M107").Value = 165.7059000000002
$ws.Range("N107").Value = -9215.142599999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29957.857
$ws.Range("J57").Value = 29957.857
$ws.Range("L57").Value = 29957.857
$ws.Range("N57").Value = -31597.857
$ws.Range("H123").Value = 13533.333
$ws.Range("J123").Value = 13533.333
$ws.Range("L123").Value = 13533.333
$ws.Range("N123").Value = -18433.333
$ws.Range("H126").Value = 4760.05
$ws.Range("I126").Value = 2447.3044
$ws.Range("J126").Value = 7889.0586
$ws.Range("K126").Value = 7341.9132
$ws.Range("L126").Value = 23667.1758
$ws.Range("M126").Value = -4871.9132
$ws.Range("N126").Value = -28607.1758
$ws.Range("H132").Value = 2165.3076
$ws.Range("I132").Value = 1765.8148
$ws.Range("J132").Value = 3064.1667
$ws.Range("K132").Value = 5297.4444
$ws.Range("L132").Value = 9192.500100000001
$ws.Range("M132").Value = -2767.4444
$ws.Range("N132").Value = -14252.5001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 30092
$ws.Range("J109").Value = 30092
$ws.Range("L109").Value = 30092
$ws.Range("N109").Value = -32866
$ws.Range("H132").Value = 1906.2188
$ws.Range("I132").Value = 1435.7727
$ws.Range("J132").Value = 2941.2
$ws.Range("K132").Value = 4307.3181
$ws.Range("L132").Value = 8823.599999999999
$ws.Range("M132").Value = -1777.3181
$ws.Range("N132").Value = -13883.6
$ws.Range("H136").Value = 1785.238
$ws.Range("I136").Value = 1618.1875
$ws.Range("J136").Value = 2319.8
$ws.Range("K136").Value = 4854.5625
$ws.Range("L136").Value = 6959.400000000001
$ws.Range("M136").Value = -2304.5625
$ws.Range("N136").Value = -12059.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 12800
$ws.Range("J47").Value = 12800
$ws.Range("L47").Value = 12800
$ws.Range("N47").Value = -13944
$ws.Range("H132").Value = 2528.6035
$ws.Range("I132").Value = 2710.077
$ws.Range("J132").Value = 2156.1052
$ws.Range("K132").Value = 8130.231000000001
$ws.Range("L132").Value = 6468.3156
$ws.Range("M132").Value = -5600.231000000001
$ws.Range("N132").Value = -11528.3156
